# Add a new weekly price record for "Hortaliza, Terminal La Palmera de La
# Serena - Berenjena" as row 251, pushing the existing rows 251:274 down to
# 252:275 (mirrors inserting a new row above the old row 251 in Excel).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 251; everything below shifts down one.
$ws.Rows.Item(251).Insert()

# Populate the newly inserted row with the new data point.
$ws.Cells.Item(251, 1).Value  = 8
$ws.Cells.Item(251, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(251, 3).Value  = "Coquimbo"
$ws.Cells.Item(251, 4).Value  = 45132
$ws.Cells.Item(251, 5).Value  = 4
$ws.Cells.Item(251, 6).Value  = 100112001
$ws.Cells.Item(251, 7).Value  = "Berenjena"
$ws.Cells.Item(251, 8).Value  = "Sin especificar"
$ws.Cells.Item(251, 9).Value  = "Primera"
$ws.Cells.Item(251, 10).Value = 400
$ws.Cells.Item(251, 11).Value = 8000
$ws.Cells.Item(251, 12).Value = 9000
$ws.Cells.Item(251, 13).Value = 8500
$ws.Cells.Item(251, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(251, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(251, 16).Value = 170
$ws.Cells.Item(251, 17).Value = 50
$ws.Cells.Item(251, 18).Value = "Hortaliza"
